$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: remove the centered alignment (<w:jc w:val="center"/>)
#    so the paragraph falls back to its style default (left).
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Alignment = 0

# ---------------------------------------------------------------------
# 2) Insert ", Oliver Szegedi," right after "My client" in the scenario
#    paragraph ("My client is an internet user..." ->
#    "My client, Oliver Szegedi, is an internet user...").
#
#    Editing text inside that paragraph's first run tends to coalesce
#    the *whole* paragraph's runs back together on save, so first drop
#    a temporary bookmark at every existing run boundary in the
#    paragraph to pin those breaks in place, make the edit (which also
#    gets its own internal run breaks pinned the same way), then remove
#    all the temporary bookmarks again -- leaving the run breaks behind
#    but no stray bookmark markers.
#
#    The final "_GoBack" bookmark is dropped right after the (already
#    present) " is an" that follows the inserted text -- Word keeps
#    only a single "_GoBack" bookmark and always relocates it to the
#    most-recently-edited spot, which is exactly how it disappears from
#    its old location further down (after "user-friendly") in the
#    target document.
# ---------------------------------------------------------------------

$scan = $d.Content
$scan.Find.Execute("My client is an internet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scan.Collapse(0)

$boundaryAnchors = @(
  "ber many different passwords. ",
  "Therefore, he wants",
  "store the different passwords.",
  " However, he",
  "efore, he asked me for advice.",
  " ",
  "ch he enthusiastically agreed."
)

$tempNames = New-Object System.Collections.ArrayList
$i = 0
foreach ($a in $boundaryAnchors) {
  $scan.Find.Execute($a, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
  $scan.Collapse(0)
  $bmName = "_pinBreak" + $i
  $d.Bookmarks.Add($bmName, $d.Range($scan.Start, $scan.Start)) | Out-Null
  $tempNames.Add($bmName) | Out-Null
  $i++
}

$rng = $d.Content
$rng.Find.Execute("My client", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$p0 = $rng.Start

$insert = ", Oliver Szegedi,"
$rng.InsertAfter($insert)

# Boundaries (relative to p0) matching each desired run break:
#   "My client" | ", Oliver " | "Szegedi" | "," | " is an" | <_GoBack> | " internet..."
$b1 = $p0 + 9   # after ", Oliver "
$b2 = $p0 + 16  # after "Szegedi"
$b3 = $p0 + $insert.Length  # after "," -> end of the inserted text

$d.Bookmarks.Add("_pinSplit0", $d.Range($p0, $p0)) | Out-Null
$d.Bookmarks.Add("_pinSplit1", $d.Range($b1, $b1)) | Out-Null
$d.Bookmarks.Add("_pinSplit2", $d.Range($b2, $b2)) | Out-Null
$d.Bookmarks.Add("_pinSplit3", $d.Range($b3, $b3)) | Out-Null

# The original text already continued "... is an internet user...";
# find that (untouched) " is an" right after our inserted text and drop
# "_GoBack" right after it.
$scan2 = $d.Range($b3, $d.Content.End)
$scan2.Find.Execute(" is an", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scan2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $d.Range($scan2.Start, $scan2.Start)) | Out-Null

# Drop every temporary pin bookmark again -- this leaves the run breaks
# they introduced in place, it only removes the bookmark markers.
foreach ($n in $tempNames) {
  $d.Bookmarks.Item($n).Delete()
}
$d.Bookmarks.Item("_pinSplit0").Delete()
$d.Bookmarks.Item("_pinSplit1").Delete()
$d.Bookmarks.Item("_pinSplit2").Delete()
$d.Bookmarks.Item("_pinSplit3").Delete()
